$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.252.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.396.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +8.17%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.27"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "681.92"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.942.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.373.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.390.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.67"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.28"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.23"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.72"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.45"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.74"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +13.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "555.18"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.02"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.03"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.604.77"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.30"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0745"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.63%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0426"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.62%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.16"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.63%  "
